# Update tracker results: fill "resultado" (G) and "profit" (H) columns
# for the rows that were recently decided.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$updates = @(
    @{ Row = 183; Resultado = "Fallo";   Profit = -1 },
    @{ Row = 184; Resultado = "Acierto"; Profit = 2.25 },
    @{ Row = 188; Resultado = "Acierto"; Profit = 1.63 },
    @{ Row = 191; Resultado = "Acierto"; Profit = 0.67 },
    @{ Row = 193; Resultado = "Fallo";   Profit = -1 },
    @{ Row = 194; Resultado = "Fallo";   Profit = -1 },
    @{ Row = 196; Resultado = "Fallo";   Profit = -1 },
    @{ Row = 197; Resultado = "Fallo";   Profit = -1 },
    @{ Row = 198; Resultado = "Acierto"; Profit = 1.2 },
    @{ Row = 200; Resultado = "Acierto"; Profit = 0.62 },
    @{ Row = 203; Resultado = "Fallo";   Profit = -1 },
    @{ Row = 204; Resultado = "Acierto"; Profit = 1 }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 7).Value = $u.Resultado
    $ws.Cells.Item($u.Row, 8).Value = $u.Profit
}
